# Fruta / hortaliza, semanal
#
# Insert 4 new weekly price rows (dated 2022-05-13, serial 44694) for the
# "Vega Monumental Concepcion - Pera" subset, right above what used to be
# the first row of this varietal block. Inserting whole rows shifts the
# remaining rows (old 339-407) down to 343-411 and carries the date
# column's number format along with them, matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 339:407 down to 343:411, leaving 339:342 free.
$ws.Rows("339:342").Insert()

# This week's data for the 4 newly-inserted rows.
$data = @(
    @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 'Fruta', 100104, 'Frutos de pepita', 100104005, 'Pera', 'Abate Fettel', 'Primera', 220, 8000, 9000, 8545, '$/caja 16 kilos empedrada', 'Región de O''Higgins', 534, 16),
    @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 'Fruta', 100104, 'Frutos de pepita', 100104005, 'Pera', 'Packham''s Triumph', 'Especial', 200, 10000, 10000, 10000, '$/caja 16 kilos empedrada', 'Región de O''Higgins', 625, 16),
    @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 'Fruta', 100104, 'Frutos de pepita', 100104005, 'Pera', 'Packham''s Triumph', 'Primera', 250, 9000, 9000, 9000, '$/caja 16 kilos empedrada', 'Región de O''Higgins', 562, 16),
    @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 'Fruta', 100104, 'Frutos de pepita', 100104005, 'Pera', 'Packham''s Triumph', 'Segunda', 150, 8000, 8000, 8000, '$/caja 16 kilos empedrada', 'Región de O''Higgins', 500, 16)
)

$arr = New-Object 'object[,]' 4,20
for ($i = 0; $i -lt 4; $i++) {
    for ($j = 0; $j -lt 20; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$ws.Range("A339:T342").Value = $arr
